# "Generate Report for Handback"
# Populates the "Latest Target File" / "Latest Handback File" columns (F/G) for
# the zh-cn and de-de worksheets, flips the status message to reflect a
# completed handback, and stamps the handback datetime.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
# BGR-packed value of RGB(0x64,0x95,0xED) == the custom "HyperLink" font color
# already used elsewhere in this workbook (FF6495ED).
$linkColor = 15570276

function Set-HandbackLink($ws, $cellAddr, $text, $url) {
    $cell = $ws.Range($cellAddr)
    $cell.Value = $text
    $ws.Hyperlinks.Add($cell, $url, "", "", $text)
    $cell.Style = "HyperLink"
    $cell.Font.Underline = $true
    $cell.Font.Color = $linkColor
}

# --- zh-cn worksheet ------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

Set-HandbackLink $wsZh "F2" "a.md" "https://github.com/OpenLocalizationTest/oltest/blob/4686b560b49f9916fd8f5d0f22769e0e84dec346/e2e/a.md"
Set-HandbackLink $wsZh "G2" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2a999a7e9110aeab370547e3502173d96774a2f4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
Set-HandbackLink $wsZh "F3" "a.md" "https://github.com/OpenLocalizationTest/oltest/blob/4686b560b49f9916fd8f5d0f22769e0e84dec346/e2e/a.md"
Set-HandbackLink $wsZh "G3" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2a999a7e9110aeab370547e3502173d96774a2f4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("H2").Value = "2016-03-22 08:30:06"
$wsZh.Range("H3").Value = "2016-03-22 08:30:06"

# --- de-de worksheet --------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

Set-HandbackLink $wsDe "F2" "a.md" "https://github.com/OpenLocalizationTest/oltest/blob/4686b560b49f9916fd8f5d0f22769e0e84dec346/e2e/a.md"
Set-HandbackLink $wsDe "G2" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1434616984613fe4989fbc5f750cf2a11537e938/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
Set-HandbackLink $wsDe "F3" "a.md" "https://github.com/OpenLocalizationTest/oltest/blob/4686b560b49f9916fd8f5d0f22769e0e84dec346/e2e/a.md"
Set-HandbackLink $wsDe "G3" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1434616984613fe4989fbc5f750cf2a11537e938/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("H2").Value = "2016-03-22 08:30:15"
$wsDe.Range("H3").Value = "2016-03-22 08:30:15"
